# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the latest scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F13").Value = 1466
    $ws.Range("F30").Value = 94
    $ws.Range("F34").Value = 239
    $ws.Range("F39").Value = 112
}
